$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22 (student #19): the attendance marks in columns G/H were blank and
# columns I/J didn't exist yet. Fill them all in with 5, matching the other
# fully-attended rows.
$ws.Range("G22").Value = 5
$ws.Range("H22").Value = 5

# I22/J22 are brand new cells in this row; give them the same formatting
# (style) used for the equivalent cells in another fully-attended row (5)
# before writing their values, so the style carries over correctly.
$ws.Range("I5").Copy()
$ws.Range("I22").PasteSpecial(-4122)
$ws.Range("I22").Value = 5

$ws.Range("J5").Copy()
$ws.Range("J22").PasteSpecial(-4122)
$ws.Range("J22").Value = 5

# The total in L22 (shared SUM formula) recalculates automatically to 40.

# Move the active selection to C22.
$ws.Range("C22").Select()
